$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Unprotect()

# Update the confidential disclaimer text in A13: date changes from 2021-05-06 to 2021-05-07
$ws.Range("A13").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution." + [char]10 + "Model holdings provided as of 2021-05-07 for illustrative purposes only and are subject to change."

# Update Weight (column D) and Percent Change (column E) values for rows 2-10
$ws.Range("D2").Value = 0.09290135810089294
$ws.Range("E2").Value = 0.01290080160320639

$ws.Range("D3").Value = 0.1067648450827276
$ws.Range("E3").Value = 0.01005747126436796

$ws.Range("D4").Value = 0.120076034552637
$ws.Range("E4").Value = 0.007701882682433547

$ws.Range("D5").Value = 0.1403956663964255
$ws.Range("E5").Value = 0.01135100809652334

$ws.Range("D6").Value = 0.1368288405657009
$ws.Range("E6").Value = 0.004407713498622678

$ws.Range("D7").Value = 0.1470809245882521
$ws.Range("E7").Value = 0.009546986147510372

$ws.Range("D8").Value = 0.1279707751296294
$ws.Range("E8").Value = 0.01384388807069215

$ws.Range("D9").Value = 0.1279815555837346
$ws.Range("E9").Value = 0.01255848614046462

$ws.Range("E10").Value = 0.01017687978448834

# Restore sheet protection (the sheet was protected before this edit)
$ws.Protect()
